$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 500582.5
$ws.Range("I4").Value = 1000115
$ws.Range("K4").Value = 1000115
$ws.Range("M4").Value = -1000001

$ws.Range("H19").Value = 1713.625
$ws.Range("I19").Value = 1018.8
$ws.Range("K19").Value = 1018.8
$ws.Range("M19").Value = -843.8

$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H39").Value = 1597.4117
$ws.Range("I39").Value = 1075
$ws.Range("K39").Value = 3225
$ws.Range("M39").Value = -2929

$ws.Range("H62").Value = 67831.836
$ws.Range("I62").Value = 73089.37
$ws.Range("K62").Value = 73089.37
$ws.Range("M62").Value = -72465.37

$ws.Range("H65").Value = 67831.836
$ws.Range("I65").Value = 73089.37
$ws.Range("K65").Value = 365446.85
$ws.Range("M65").Value = -362326.85

$ws.Range("H129").Value = 5491.2915
$ws.Range("I129").Value = 1669
$ws.Range("J129").Value = 10842.5
$ws.Range("K129").Value = 5007
$ws.Range("L129").Value = 32527.5
$ws.Range("M129").Value = -7
$ws.Range("N129").Value = -42527.5

$ws.Range("H132").Value = 7037.4375
$ws.Range("I132").Value = 5221.4634
$ws.Range("K132").Value = 15664.3902
$ws.Range("M132").Value = -13134.3902

$ws.Range("H137").Value = 3094.6667
$ws.Range("I137").Value = 2634.077
$ws.Range("J137").Value = 3522.3572
$ws.Range("K137").Value = 7902.231000000001
$ws.Range("L137").Value = 10567.0716
$ws.Range("M137").Value = -5352.231000000001
$ws.Range("N137").Value = -15667.0716

$ws.Range("H141").Value = 2746.44
$ws.Range("I141").Value = 2746.44
$ws.Range("K141").Value = 8239.32
$ws.Range("M141").Value = -3059.32

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 9999
$ws.Range("J6").Value = 9999
$ws.Range("L6").Value = 9999
$ws.Range("N6").Value = -10345

$ws.Range("H45").Value = 1002248.4
$ws.Range("I45").Value = 1429569.4
$ws.Range("K45").Value = 1429569.4
$ws.Range("M45").Value = -1429192.4

$ws.Range("H61").Value = 6486.1763
$ws.Range("I61").Value = 6289
$ws.Range("K61").Value = 6289
$ws.Range("M61").Value = -6077

$ws.Range("H74").Value = 48752.863
$ws.Range("I74").Value = 58592.445
$ws.Range("K74").Value = 58592.445
$ws.Range("M74").Value = -57718.445

$ws.Range("H77").Value = 48752.863
$ws.Range("I77").Value = 58592.445
$ws.Range("K77").Value = 292962.225
$ws.Range("M77").Value = -288594.225

$ws.Range("H95").Value = 66671
$ws.Range("J95").Value = 66671
$ws.Range("L95").Value = 66671
$ws.Range("N95").Value = -72163

$ws.Range("H122").Value = 4356.2856
$ws.Range("I122").Value = 4415.6665
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 13246.9995
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -10796.9995
$ws.Range("N122").Value = -16900

$ws.Range("H132").Value = 25852.295
$ws.Range("I132").Value = 29990.648
$ws.Range("J132").Value = 3978.1428
$ws.Range("K132").Value = 89971.944
$ws.Range("L132").Value = 11934.4284
$ws.Range("M132").Value = -87441.944
$ws.Range("N132").Value = -16994.4284

$ws.Range("H136").Value = 6486.1763
$ws.Range("I136").Value = 6289
$ws.Range("K136").Value = 18867
$ws.Range("M136").Value = -16317

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2698.4666
$ws.Range("I86").Value = 1969.125
$ws.Range("K86").Value = 1969.125
$ws.Range("M86").Value = -846.125

$ws.Range("H89").Value = 2698.4666
$ws.Range("I89").Value = 1969.125
$ws.Range("K89").Value = 9845.625
$ws.Range("M89").Value = -4229.625

$ws.Range("H134").Value = 5750
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 5750
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 17250
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -22320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 30000
$ws.Range("I17").Value = 30000
$ws.Range("K17").Value = 30000
$ws.Range("M17").Value = -29826

$ws.Range("H43").Value = 45656
$ws.Range("J43").Value = 45656
$ws.Range("L43").Value = 45656
$ws.Range("N43").Value = -46024

$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H101").Value = 45656
$ws.Range("J101").Value = 45656
$ws.Range("L101").Value = 45656
$ws.Range("N101").Value = -52146

$ws.Range("H132").Value = 4522.154
$ws.Range("I132").Value = 4697.5
$ws.Range("K132").Value = 14092.5
$ws.Range("M132").Value = -11562.5

$ws.Range("H134").Value = 36465.645
$ws.Range("I134").Value = 43945.6
$ws.Range("J134").Value = 5299.1665
$ws.Range("K134").Value = 131836.8
$ws.Range("L134").Value = 15897.4995
$ws.Range("M134").Value = -129301.8
$ws.Range("N134").Value = -20967.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 43.3125
$ws.Range("J38").Value = 38.857143
$ws.Range("L38").Value = 116.571429
$ws.Range("N38").Value = -810.571429

$ws.Range("H44").Value = 1347.5
$ws.Range("I44").Value = 796.6667
$ws.Range("J44").Value = 3000
$ws.Range("K44").Value = 2390.0001
$ws.Range("L44").Value = 9000
$ws.Range("M44").Value = -1992.0001
$ws.Range("N44").Value = -9796

$ws.Range("H122").Value = 2497.5
$ws.Range("J122").Value = 2497.5
$ws.Range("L122").Value = 22477.5
$ws.Range("N122").Value = -27377.5

$ws.Range("H133").Value = 8785.799999999999
$ws.Range("I133").Value = 8785.799999999999
$ws.Range("K133").Value = 26357.4
$ws.Range("M133").Value = -21297.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2833.111
$ws.Range("I102").Value = 2700
$ws.Range("K102").Value = 2700
$ws.Range("M102").Value = -1078

$ws.Range("H122").Value = 3538.6
$ws.Range("I122").Value = 3376.2222
$ws.Range("K122").Value = 10128.6666
$ws.Range("M122").Value = -7678.6666

$ws.Range("H132").Value = 152847.92
$ws.Range("I132").Value = 94706
$ws.Range("K132").Value = 284118
$ws.Range("M132").Value = -281588

$ws.Range("H138").Value = 87449
$ws.Range("J138").Value = 87449
$ws.Range("L138").Value = 87449
$ws.Range("N138").Value = -97729

$ws.Range("H139").Value = 73594.664
$ws.Range("J139").Value = 73594.664
$ws.Range("L139").Value = 73594.664
$ws.Range("N139").Value = -83874.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 32223.191
$ws.Range("I132").Value = 37734.2
$ws.Range("K132").Value = 113202.6
$ws.Range("M132").Value = -110672.6

$ws.Range("H136").Value = 5904.3887
$ws.Range("I136").Value = 5468.875
$ws.Range("K136").Value = 16406.625
$ws.Range("M136").Value = -13856.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 593.9048
$ws.Range("I107").Value = 444.16666
$ws.Range("J107").Value = 1492.3334
$ws.Range("K107").Value = 1332.49998
$ws.Range("L107").Value = 4477.0002
$ws.Range("M107").Value = 587.5000199999999
$ws.Range("N107").Value = -8317.0002

$ws.Range("H113").Value = 1173
$ws.Range("I113").Value = 1108.0312
$ws.Range("J113").Value = 1311.6
$ws.Range("K113").Value = 3324.0936
$ws.Range("L113").Value = 3934.8
$ws.Range("M113").Value = -1154.0936
$ws.Range("N113").Value = -8274.799999999999

$ws.Range("H132").Value = 41390.574
$ws.Range("I132").Value = 43138.75
$ws.Range("J132").Value = 31600.8
$ws.Range("K132").Value = 129416.25
$ws.Range("L132").Value = 94802.39999999999
$ws.Range("M132").Value = -126886.25
$ws.Range("N132").Value = -99862.39999999999

$ws.Range("H136").Value = 30495432
$ws.Range("I136").Value = 3493125.8
$ws.Range("J136").Value = 125003500
$ws.Range("K136").Value = 10479377.4
$ws.Range("L136").Value = 375010500
$ws.Range("M136").Value = -10476827.4
$ws.Range("N136").Value = -375015600
